$wb = $excel.ActiveWorkbook

# --- 1. Sheet2 ("Equipos"): update standings table abbreviations (done first so the
#        shared-string table picks up insertion order QAR,GAB,ZIR,NEF,SAB,SAM,SUM,SEB) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("C2").Value = "QAR"
$ws2.Range("C3").Value = "GAB"
$ws2.Range("C4").Value = "ZIR"
$ws2.Range("C5").Value = "NEF"
$ws2.Range("C6").Value = "SAB"
$ws2.Range("C7").Value = "SAM"
$ws2.Range("C8").Value = "SUM"
$ws2.Range("C9").Value = "SEB"

# --- 2. Sheet1 ("Resultados"): set every match date to 2022-05-21 (serial 44702) ---
$ws1 = $wb.Worksheets.Item(1)
$blocks = @("B3:B6","B8:B11","B13:B16","B18:B21","B23:B26","B28:B31","B33:B36","B38:B41","B43:B46","B48:B51","B53:B56","B58:B61","B63:B66","B68:B71")
$ws1.Range("B3").Copy()
foreach ($b in $blocks) {
    $rng = $ws1.Range($b)
    $rng.PasteSpecial(-4122)
    $rng.Value2 = 44702
}

# --- 3. Sheet1: replace Local/Visita team names with their abbreviations ---
$pairs = @(
    @("C3","GAB"),
    @("D3","SAB"),
    @("C4","ZIR"),
    @("D4","SAM"),
    @("C5","SUM"),
    @("D5","SEB"),
    @("C6","QAR"),
    @("D6","NEF"),
    @("C8","SAM"),
    @("D8","GAB"),
    @("C9","NEF"),
    @("D9","ZIR"),
    @("C10","SAB"),
    @("D10","SUM"),
    @("C11","SEB"),
    @("D11","QAR"),
    @("C13","SEB"),
    @("D13","SAB"),
    @("C14","QAR"),
    @("D14","ZIR"),
    @("C15","SUM"),
    @("D15","SAM"),
    @("C16","GAB"),
    @("D16","NEF"),
    @("C18","SAM"),
    @("D18","SEB"),
    @("C19","QAR"),
    @("D19","SAB"),
    @("C20","ZIR"),
    @("D20","GAB"),
    @("C21","SUM"),
    @("D21","NEF"),
    @("C23","GAB"),
    @("D23","QAR"),
    @("C24","SAB"),
    @("D24","SAM"),
    @("C25","SEB"),
    @("D25","NEF"),
    @("C26","SUM"),
    @("D26","ZIR"),
    @("C28","GAB"),
    @("D28","SUM"),
    @("C29","QAR"),
    @("D29","SAM"),
    @("C30","ZIR"),
    @("D30","SEB"),
    @("C31","NEF"),
    @("D31","SAB"),
    @("C33","SAM"),
    @("D33","NEF"),
    @("C34","SAB"),
    @("D34","ZIR"),
    @("C35","SUM"),
    @("D35","QAR"),
    @("C36","SEB"),
    @("D36","GAB"),
    @("C38","QAR"),
    @("D38","SEB"),
    @("C39","SUM"),
    @("D39","SAB"),
    @("C40","GAB"),
    @("D40","SAM"),
    @("C41","NEF"),
    @("D41","ZIR"),
    @("C43","SAM"),
    @("D43","SUM"),
    @("C44","NEF"),
    @("D44","GAB"),
    @("C45","SAB"),
    @("D45","SEB"),
    @("C46","ZIR"),
    @("D46","QAR"),
    @("C48","SEB"),
    @("D48","SAM"),
    @("C49","SAB"),
    @("D49","QAR"),
    @("C50","SUM"),
    @("D50","NEF"),
    @("C51","GAB"),
    @("D51","ZIR"),
    @("C53","ZIR"),
    @("D53","SUM"),
    @("C54","NEF"),
    @("D54","SEB"),
    @("C55","SAM"),
    @("D55","SAB"),
    @("C56","QAR"),
    @("D56","GAB"),
    @("C58","SAB"),
    @("D58","NEF"),
    @("C59","SUM"),
    @("D59","GAB"),
    @("C60","SEB"),
    @("D60","ZIR"),
    @("C61","SAM"),
    @("D61","QAR"),
    @("C63","GAB"),
    @("D63","SEB"),
    @("C64","NEF"),
    @("D64","SAM"),
    @("C65","ZIR"),
    @("D65","SAB"),
    @("C66","QAR"),
    @("D66","SUM"),
    @("C68","SAM"),
    @("D68","ZIR"),
    @("C69","NEF"),
    @("D69","QAR"),
    @("C70","SAB"),
    @("D70","GAB"),
    @("C71","SEB"),
    @("D71","SUM")
)
foreach ($p in $pairs) {
    $ws1.Range($p[0]).Value = $p[1]
}

# --- 4. Restore view state (zoom + selection) on both sheets ---
$ws2.Activate()
$ws2.Range("B9").Select() | Out-Null

$ws1.Activate()
$excel.ActiveWindow.Zoom = 80
$ws1.Range("B74").Select() | Out-Null
